# Apply "Add data for 2021-10-24" update to the carjacking-by-neighborhood
# workbook: rename the "through October 15" sheet/title/header to
# "through October 16", and bump a handful of neighborhood/month counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet tab and title references ---------------------------------
$ws.Name = "Through 2021-10-16"
$ws.Range("B1").Value = "October 2021 (through October 16)"

# --- Update / add individual count cells ------------------------------------
# Row 2 - Garfield Park
$ws.Range("B2").Value = 9
$ws.Range("L2").Value = 11
$ws.Range("AF2").Value = 5

# Row 4 - North Lawndale
$ws.Range("B4").Value = 8
$ws.Range("L4").Value = 8

# Row 10 - Roseland
$ws.Range("L10").Value = 4

# Row 13 - South Shore
$ws.Range("B13").Value = 4

# Row 15 - West Town (new value in previously empty cell)
$ws.Range("L15").Value = 1

# Row 22 - Grand Boulevard
$ws.Range("AZ22").Value = 2

# Row 24 - Chinatown (new value in previously empty cell)
$ws.Range("V24").Value = 1

# Row 38 - Edgewater
$ws.Range("B38").Value = 2

# Row 45 - West Lawn (new value in previously empty cell)
$ws.Range("AF45").Value = 1

# Row 47 - West Ridge (new value in previously empty cell)
$ws.Range("AP47").Value = 1

# Row 64 - Brighton Park
$ws.Range("AZ64").Value = 2

# Row 73 - Galewood
$ws.Range("AF73").Value = 2

# Row 75 - Greektown (new value in previously empty cell)
$ws.Range("L75").Value = 1
